$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 3.1
$ws.Range("K2").Value = 1.95
$ws.Range("M2").Value = 1.05
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 2.62
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53
$ws.Range("X2").Value = 10
$ws.Range("Y2").Value = 10
$ws.Range("AF2").Value = 67
$ws.Range("AI2").Value = 15
$ws.Range("AP2").Value = 29
